$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'274.97"
$ws.Range("D2").ClearFormats()
$ws.Range("D4").Value = "'6.298"
$ws.Range("D4").ClearFormats()
$ws.Range("D5").Value = "'0.06267"
$ws.Range("D5").ClearFormats()
$ws.Range("D6").Value = "'3.663"
$ws.Range("D6").ClearFormats()
$ws.Range("D7").Value = "'6.679"
$ws.Range("D7").ClearFormats()
$ws.Range("D8").Value = "'1.361"
$ws.Range("D8").ClearFormats()
$ws.Range("D9").Value = "'0.8313"
$ws.Range("D9").ClearFormats()
$ws.Range("D10").Value = "'0.01376"
$ws.Range("D10").ClearFormats()
$ws.Range("D11").Value = "'0.1631"
$ws.Range("D11").ClearFormats()
$ws.Range("D12").Value = "'0.08369"
$ws.Range("D12").ClearFormats()
$ws.Range("D13").Value = "'0.03442"
$ws.Range("D13").ClearFormats()
$ws.Range("D14").Value = "'0.03104"
$ws.Range("D14").ClearFormats()
$ws.Range("D15").Value = "'0.09315"
$ws.Range("D15").ClearFormats()
$ws.Range("D16").Value = "'3.888"
$ws.Range("D16").ClearFormats()
$ws.Range("D17").Value = "'0.001638"
$ws.Range("D17").ClearFormats()
$ws.Range("D18").Value = "'0.04766"
$ws.Range("D18").ClearFormats()
$ws.Range("D19").Value = "'0.006369"
$ws.Range("D19").ClearFormats()
$ws.Range("D20").Value = "'0.005694"
$ws.Range("D20").ClearFormats()
$ws.Range("D22").Value = "'0.0001500"
$ws.Range("D22").ClearFormats()
$ws.Range("D23").Value = "'3.715"
$ws.Range("D23").ClearFormats()
$ws.Range("D25").Value = "'0.3341"
$ws.Range("D25").ClearFormats()
$ws.Range("D26").Value = "'0.1240"
$ws.Range("D26").ClearFormats()
$ws.Range("D40").Value = "'0.04706"
$ws.Range("D40").ClearFormats()
$ws.Range("D41").Value = "'0.007059"
$ws.Range("D41").ClearFormats()
$ws.Range("D43").Value = "'0.003300"
$ws.Range("D43").ClearFormats()
$ws.Range("D44").Value = "'0.01224"
$ws.Range("D44").ClearFormats()
$ws.Range("D45").Value = "'0.00006257"
$ws.Range("D45").ClearFormats()
$ws.Range("D47").Value = "'0.7696"
$ws.Range("D47").ClearFormats()
$ws.Range("D48").Value = "'0.02991"
$ws.Range("D48").ClearFormats()
$ws.Range("D49").Value = "'0.00002200"
$ws.Range("D49").ClearFormats()

Write-Host "Updated symbol list prices"
